# Auto-generated Excel COM-interop PowerShell script
# Applies the edits described by the unified diff to before.xlsx
# Target: insert 2 new rows for Thailand matches before the old row 13,
# shifting the previous rows 13-15 down to 15-17, then update odds values
# that changed on the unaffected rows 5 and 10 plus on rows 5/10/15/16/17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: odds updates ---
$ws.Range("G5").Value = 1.7
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("X5").Value = 7
$ws.Range("Z5").Value = 13
$ws.Range("AJ5").Value = 19
$ws.Range("AN5").Value = 3.5
$ws.Range("AZ5").Value = 126

# --- Row 10: odds update ---
$ws.Range("S10").Value = 1.11

# --- Insert two new rows at 13 (old rows 13,14,15 shift down to 15,16,17) ---
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# --- Row 13 (new): THAILAND - THAI LEAGUE 1, Lamphun Warrior vs Prachuap ---
$ws.Range("A13").Value = "AwSmRIc5"
$ws.Range("B13").Value = "24/11/2024"
$ws.Range("C13").Value = "09:00"
$ws.Range("D13").Value = "THAILAND - THAI LEAGUE 1"
$ws.Range("E13").Value = "Lamphun Warrior"
$ws.Range("F13").Value = "Prachuap"
$ws.Range("G13").Value = 2.27
$ws.Range("H13").Value = 3.35
$ws.Range("I13").Value = 2.92
$ws.Range("J13").Value = 2.92
$ws.Range("K13").Value = 2.1
$ws.Range("L13").Value = 3.45
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 7.4
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.3
$ws.Range("Q13").Value = 1.87
$ws.Range("R13").Value = 1.87
$ws.Range("S13").Value = 1.4
$ws.Range("T13").Value = 2.7
$ws.Range("U13").Value = 1.7
$ws.Range("V13").Value = 2.05
$ws.Range("W13").Value = 8.25
$ws.Range("X13").Value = 11.25
$ws.Range("Y13").Value = 9
$ws.Range("Z13").Value = 23
$ws.Range("AA13").Value = 18
$ws.Range("AB13").Value = 27
$ws.Range("AC13").Value = 7.4
$ws.Range("AD13").Value = 6.4
$ws.Range("AE13").Value = 13
$ws.Range("AF13").Value = 55
$ws.Range("AG13").Value = 400
$ws.Range("AH13").Value = 10
$ws.Range("AI13").Value = 16
$ws.Range("AJ13").Value = 10.25
$ws.Range("AK13").Value = 37
$ws.Range("AL13").Value = 24
$ws.Range("AM13").Value = 30
$ws.Range("AN13").Value = 4.25
$ws.Range("AO13").Value = 12.5
$ws.Range("AP13").Value = 21
$ws.Range("AQ13").Value = 50
$ws.Range("AR13").Value = 90
$ws.Range("AS13").Value = 250
$ws.Range("AT13").Value = 2.7
$ws.Range("AU13").Value = 7
$ws.Range("AV13").Value = 60
$ws.Range("AW13").Value = 4.9
$ws.Range("AX13").Value = 16
$ws.Range("AY13").Value = 22
$ws.Range("AZ13").Value = 75
$ws.Range("BA13").Value = 100
$ws.Range("BB13").Value = 300

# --- Row 14 (new): THAILAND - THAI LEAGUE 1, Nakhon Ratchasima FC vs Pathum United ---
$ws.Range("A14").Value = "6qpMCJrC"
$ws.Range("B14").Value = "24/11/2024"
$ws.Range("C14").Value = "09:00"
$ws.Range("D14").Value = "THAILAND - THAI LEAGUE 1"
$ws.Range("E14").Value = "Nakhon Ratchasima FC"
$ws.Range("F14").Value = "Pathum United"
$ws.Range("G14").Value = 5.1
$ws.Range("H14").Value = 4.05
$ws.Range("I14").Value = 1.57
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 2.3
$ws.Range("L14").Value = 2.1
$ws.Range("M14").Value = 1.04
$ws.Range("N14").Value = 8.25
$ws.Range("O14").Value = 1.22
$ws.Range("P14").Value = 3.85
$ws.Range("Q14").Value = 1.7
$ws.Range("R14").Value = 2.07
$ws.Range("S14").Value = 1.33
$ws.Range("T14").Value = 3.05
$ws.Range("U14").Value = 1.75
$ws.Range("V14").Value = 1.98
$ws.Range("W14").Value = 16
$ws.Range("X14").Value = 32
$ws.Range("Y14").Value = 16
$ws.Range("Z14").Value = 90
$ws.Range("AA14").Value = 45
$ws.Range("AB14").Value = 45
$ws.Range("AC14").Value = 8.25
$ws.Range("AD14").Value = 7.9
$ws.Range("AE14").Value = 15.5
$ws.Range("AF14").Value = 65
$ws.Range("AG14").Value = 450
$ws.Range("AH14").Value = 7.7
$ws.Range("AI14").Value = 7.9
$ws.Range("AJ14").Value = 8
$ws.Range("AK14").Value = 11.75
$ws.Range("AL14").Value = 12
$ws.Range("AM14").Value = 23
$ws.Range("AN14").Value = 6.8
$ws.Range("AO14").Value = 28
$ws.Range("AP14").Value = 30
$ws.Range("AQ14").Value = 150
$ws.Range("AR14").Value = 175
$ws.Range("AS14").Value = 350
$ws.Range("AT14").Value = 3.05
$ws.Range("AU14").Value = 7.4
$ws.Range("AV14").Value = 60
$ws.Range("AW14").Value = 3.5
$ws.Range("AX14").Value = 7.5
$ws.Range("AY14").Value = 16
$ws.Range("AZ14").Value = 23
$ws.Range("BA14").Value = 50
$ws.Range("BB14").Value = 200

# --- Row 15 (shifted): Sivasspor vs Kasimpasa, TURKEY - SUPER LIG - odds updates only ---
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 3.3
$ws.Range("I15").Value = 2.35
$ws.Range("J15").Value = 3.5
$ws.Range("L15").Value = 3
$ws.Range("O15").Value = 1.29
$ws.Range("P15").Value = 3.5
$ws.Range("Q15").Value = 1.95
$ws.Range("R15").Value = 1.9
$ws.Range("S15").Value = 1.4
$ws.Range("T15").Value = 2.75
$ws.Range("U15").Value = 1.73
$ws.Range("V15").Value = 2
$ws.Range("W15").Value = 10
$ws.Range("X15").Value = 15
$ws.Range("Z15").Value = 29
$ws.Range("AC15").Value = 10
$ws.Range("AG15").Value = 201
$ws.Range("AH15").Value = 8.5
$ws.Range("AI15").Value = 12
$ws.Range("AJ15").Value = 9.5
$ws.Range("AK15").Value = 23
$ws.Range("AL15").Value = 19
$ws.Range("AM15").Value = 26
$ws.Range("AT15").Value = 2.75
$ws.Range("AU15").Value = 8
$ws.Range("AW15").Value = 4.5
$ws.Range("AX15").Value = 13
$ws.Range("AY15").Value = 23
$ws.Range("BA15").Value = 67
$ws.Range("BB15").Value = 151
$ws.Range("BC15").Value = 276

# --- Row 16 (shifted): Erzurumspor vs Karagumruk, TURKEY - 1. LIG - odds updates only ---
$ws.Range("G16").Value = 2.05
$ws.Range("H16").Value = 3.4
$ws.Range("I16").Value = 3.4
$ws.Range("J16").Value = 2.75
$ws.Range("K16").Value = 2.1
$ws.Range("L16").Value = 4.33
$ws.Range("Q16").Value = 2.08
$ws.Range("R16").Value = 1.73
$ws.Range("X16").Value = 9.5
$ws.Range("Y16").Value = 9
$ws.Range("Z16").Value = 19
$ws.Range("AA16").Value = 17
$ws.Range("AH16").Value = 9.5
$ws.Range("AI16").Value = 17
$ws.Range("AJ16").Value = 13
$ws.Range("AK16").Value = 41
$ws.Range("AL16").Value = 29
$ws.Range("AM16").Value = 41
$ws.Range("AN16").Value = 4
$ws.Range("AO16").Value = 12
$ws.Range("AS16").Value = 151
$ws.Range("AW16").Value = 5.5
$ws.Range("AZ16").Value = 67

# --- Row 17 (shifted): Obolon vs Kryvbas, UKRAINE - PREMIER LEAGUE - odds updates only ---
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 1.88
$ws.Range("J17").Value = 4.6
$ws.Range("K17").Value = 1.95
$ws.Range("L17").Value = 2.55
$ws.Range("M17").Value = 1.09
$ws.Range("N17").Value = 7.5
$ws.Range("O17").Value = 1.44
$ws.Range("P17").Value = 2.4
$ws.Range("Q17").Value = 2.27
$ws.Range("R17").Value = 1.5
$ws.Range("S17").Value = 1.5
$ws.Range("T17").Value = 2.27
$ws.Range("U17").Value = 2.07
$ws.Range("V17").Value = 1.6
$ws.Range("W17").Value = 8.75
$ws.Range("X17").Value = 20
$ws.Range("Y17").Value = 14.5
$ws.Range("Z17").Value = 65
$ws.Range("AA17").Value = 50
$ws.Range("AB17").Value = 65
$ws.Range("AC17").Value = 7
$ws.Range("AD17").Value = 6.4
$ws.Range("AE17").Value = 20
$ws.Range("AF17").Value = 120
$ws.Range("AH17").Value = 5.4
$ws.Range("AI17").Value = 7.6
$ws.Range("AJ17").Value = 9
$ws.Range("AK17").Value = 15.5
$ws.Range("AL17").Value = 18.5
$ws.Range("AN17").Value = 5.6
$ws.Range("AO17").Value = 24
$ws.Range("AP17").Value = 37
$ws.Range("AQ17").Value = 150
$ws.Range("AT17").Value = 2.25
$ws.Range("AU17").Value = 8.25
$ws.Range("AW17").Value = 3.55
$ws.Range("AX17").Value = 9.75
$ws.Range("AY17").Value = 23
$ws.Range("AZ17").Value = 40
$ws.Range("BA17").Value = 90

Write-Output "edit complete"
